$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New aggregated data (id, value) replacing the old rows 2-31 with rows 2-17
$data = @(
    @("id_DK_Central_BH_Biogas", 30.059),
    @("id_DK_Central_IndustryH_Biogas", 33.02200000000001),
    @("id_DK_Central_BH_Biomass", 2663.333),
    @("id_DK_Central_IndustryH_Biomass", 64.684),
    @("id_DK_Central_BH_Natgas", 5490.355),
    @("id_DK_Central_IndustryH_Natgas", 8.619),
    @("id_DK_Central_BH_Oil", 4518.727),
    @("id_DK_Central_IndustryH_Oil", 53.102),
    @("id_DK_Central_BH_Waste", 99.54900000000001),
    @("id_DK_Central_EP", 10.8),
    @("id_DK_Central_GT", 8.300000000000001),
    @("id_DK_Central_HPstandard", 566.432),
    @("id_DK_Central_HPsurplusheat", 79.771),
    @("id_DK_Central_IH", 1567.35),
    @("id_DK_Central_IndustryH", 362.12),
    @("id_DK_Central_SH", 1107.157)
)

# Clear the old rows below the data we are about to write (rows 2 through old last row 31)
$ws.Range("A2:B31").Clear()

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
